$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column S (year 2022) values, mirroring column R's styling
$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 27.292394741221504
$ws.Range("S6").Value = 36.613942589338023
$ws.Range("S7").Value = 14.18691257315127
$ws.Range("S8").Value = 55.377118174770182
$ws.Range("S9").Value = 42.247570764681029
$ws.Range("S10").Value = 30.18817294468856
$ws.Range("S11").Value = 97.03085581214826
$ws.Range("S12").Value = 25.2
$ws.Range("S13").Value = 21.849963583394029
$ws.Range("S14").Value = "-"

# Copy formatting from column R so column S matches the styling used by the rest of the table
$ws.Range("R3:R14").Copy()
$ws.Range("S3:S14").PasteSpecial(-4122)  # xlPasteFormats

# Update selection to reflect the post-edit active cell
$ws.Range("S16").Select()
